# Delete the last slide ("Future Plans of the Product contd..", slide 13)
# from the presentation, matching the target diff which drops the
# <p:sldId id="268" r:id="rId14"/> entry and the slide13.xml part.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
